# Reproduces the commit "modified complete application fro immidart enterprise"
# on TestData/ImmidartTestData.xlsx.
#
# Net semantic effect of the OOXML diff:
#  - The LoginTestData sheet's credential row is updated: the e-mail value
#    in A2 changes from "admin@immidart.com" to "admin@immidartqa.com"
#    (the password in B2 stays "password" - the shared-string table is
#    simply re-ordered/re-deduplicated by Excel as a natural side effect
#    of rewriting the cell content).
#  - The saved active selection on the sheet moves from D14 to B4.
#  - The saved window width for the workbook view changes slightly
#    (12345 -> 12030); this is pure UI chrome state, so it is still set
#    through the documented object model for completeness even though it
#    is cosmetic.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the login e-mail used as test data (QA environment address).
$ws.Range("A2").Value = "admin@immidartqa.com"

# Move/restore the active selection to B4, as recorded in the saved view.
$ws.Range("B4").Select() | Out-Null

# Best-effort: restore the saved window size (cosmetic workbook view state).
$win = $excel.ActiveWindow
$win.Width = 12030
$win.Height = 6225
